# Weekly update: insert a new daily price record for "Melón" / "Tuna" /
# "Primera" at row 63, pushing the existing rows 63-87 down to 64-88.
#
# Resulting used range grows from A1:R87 to A1:R88.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new blank row above the current row 63; Excel shifts every
# row from 63 downward (old 63 -> 64, ..., old 87 -> 88) and the new blank
# row inherits formatting (e.g. the date-format style on column D) from the
# surrounding rows.
$ws.Rows(63).Insert()

# Populate the freshly inserted row 63 with the new record.
$ws.Cells.Item(63, 1).Value  = 1
$ws.Cells.Item(63, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(63, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(63, 4).Value  = 45229
$ws.Cells.Item(63, 5).Value  = 15
$ws.Cells.Item(63, 6).Value  = 100112027
$ws.Cells.Item(63, 7).Value  = "Melón"
$ws.Cells.Item(63, 8).Value  = "Tuna"
$ws.Cells.Item(63, 9).Value  = "Primera"
$ws.Cells.Item(63, 10).Value = 50
$ws.Cells.Item(63, 11).Value = 19000
$ws.Cells.Item(63, 12).Value = 20000
$ws.Cells.Item(63, 13).Value = 19500
$ws.Cells.Item(63, 14).Value = "$/caja 18 unidades"
$ws.Cells.Item(63, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(63, 16).Value = 1083
$ws.Cells.Item(63, 17).Value = 18
$ws.Cells.Item(63, 18).Value = "Hortaliza"
